# Locate the paragraph that ends with "打完大大无大无多爱我的安慰打我打"
# (the last paragraph in the document, which currently carries the
# _GoBack bookmark), then split it into three paragraphs:
#   1. the original text, with its paragraph-mark font hint switched
#      from "default" to "eastAsia"
#   2. a new, empty paragraph (hint "eastAsia")
#   3. a new paragraph containing "12312313213212323123113123132311323",
#      which now carries the trailing _GoBack bookmark
$d = $word.ActiveDocument

$needle = "打完大大无大无多爱我的安慰打我打"
$newLine = "12312313213212323123113123132311323"

$hit = $d.Content
$found = $hit.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not find target paragraph text"
}

# Locate the Word Paragraph object that fully contains the found hit so
# we also pick up its trailing paragraph mark (and bookmark).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $hit.Start -and $p.Range.End -ge $hit.End) {
        $target = $p.Range
    }
}
if ($target -eq $null) {
    throw "could not resolve containing paragraph"
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
        '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr>' + `
            '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>' + $needle + '</w:t></w:r>' + `
        '</w:p>' + `
        '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p>' + `
        '<w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr>' + `
            '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>' + $newLine + '</w:t></w:r>' + `
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
        '</w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
